# Populate the "Catcher's View" pitch-by-pitch visual tables for both
# at-bats shown on this hitter report (rows 10-13 and rows 19-23), and
# normalize the pitch-mix order to match the order pitches actually
# occurred in each at-bat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat #1 (rows 10-13) ---
$ws.Range("F10").Value = "CH"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Strike"
$ws.Range("M10").Value = "75.42 MPH"

$ws.Range("F11").Value = "FB"
$ws.Range("G11").Value = "Swing"
$ws.Range("H11").Value = "Foul"

$ws.Range("F12").Value = "SL"
$ws.Range("G12").Value = "Swing"
$ws.Range("H12").Value = "Foul"
$ws.Range("M12").Value = "14.57°"

$ws.Range("F13").Value = "SL"
$ws.Range("G13").Value = "Swing"
$ws.Range("H13").Value = "In Play"

$ws.Range("J17").Value = "CH,FB,SL"

# --- At-bat #2 (rows 19-23) ---
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Ball"
$ws.Range("M19").Value = "61.48 MPH"

$ws.Range("F20").Value = "CH"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Strike"

$ws.Range("F21").Value = "CH"
$ws.Range("G21").Value = "Take"
$ws.Range("H21").Value = "Ball"
$ws.Range("M21").Value = "14.14°"

$ws.Range("F22").Value = "FB"
$ws.Range("G22").Value = "Take"
$ws.Range("H22").Value = "Strike"

$ws.Range("F23").Value = "CH"
$ws.Range("G23").Value = "Swing"
$ws.Range("H23").Value = "In Play"

$ws.Range("J26").Value = "CH,FB,SL"
